$d = $word.ActiveDocument

# 1. Fix "index" -> "índex" (accent) in the methodology section.
$d.Content.Find.Execute("index", $true, $true, $false, $false, $false, $true, 1, $false, "índex", 2) | Out-Null

# 2. Insert the three new narrative paragraphs (plus two trailing blank
#    paragraphs) right after the paragraph ending in
#    "...nos aeroportos mais importantes." at the close of the
#    Metodologia section, before the blank paragraphs that precede the
#    "Análise de resultados" heading.
$anchor = $d.Content
$found = $anchor.Find.Execute("nos aeroportos mais importantes.", $true)
if (-not $found) {
    throw "anchor paragraph not found"
}
$anchor.Collapse(0)

$para1 = "Depois começámos a fazer a análise avançada, nesta secção o objetivo é juntar os dois datasets e criar formas de analisar o que temos. Primeiro, analisamos a taxa de cancelamento por aeroporto em relação ao preço médio dos bilhetes por destino. A taxa de cancelamento é calculada para cada aeroporto, e os dados são combinados com as informações sobre o preço médio das passagens. Isso permite verificar se há alguma correlação entre os preços e a taxa de cancelamento nos diferentes aeroportos, ajudando a identificar se aeroportos com maior taxa de cancelamento também têm preços mais altos ou vice-versa. A seguir, fazemos uma comparação entre o número de voos por origem e a taxa de cancelamento por aeroporto. O número de voos de cada aeroporto de origem é calculado, e esses dados são combinados com as taxas de cancelamento. A análise procura entender se aeroportos com maior volume de voos têm uma taxa de cancelamento maior ou menor, criando conclusões sobre como a quantidade de voos pode influenciar o desempenho em termos de cancelamentos. Por último, a média de atraso por aeroporto é comparada com a distância média dos voos por origem. A distância média dos voos de cada aeroporto é calculada e combinada com os dados de atraso médio. Esta análise visa identificar se voos mais longos tendem a ter maiores atrasos, permitindo explorar a relação entre a distância dos voos e o tempo de atraso nos aeroportos."

$para2 = "Para concluir a secção da metodologia temos o script para exemplificar a ingestão de dados em batch e streaming. Primeiro, tal como no notebook, o script cria uma sessão de spark, configurando também a ligação com o MongoDB onde guardamos os dados. No processamento em batch carregamos os dados de voos e atrasos de uma só vez, criando relatórios. Nos atrasos começamos por agrupar os dados por aeroporto e calculamos o total de voos atrasados, tal como a média de atrasos para cada. No caso dos voos, calculamos os rendimentos totais por companhia aérea, multiplicando o preço dos bilhetes pelo número de bilhetes vendidos, acabando com o preço médio dos bilhetes. No fim mostramos parte dos relatórios na consola e guardamos isso em novas coleções no MongoDB. "

$para3 = "Passando agora para o processamento em streaming, o objetivo desta parte seria no fim de cada dia juntar aos resumos diários os dados mais recentes, ao contrário do batch que lê sempre todos. Neste caso, como os dados têm uma janela temporal específica simplesmente filtramos os dados para aqueles que estão a ser analisados e depois criamos relatórios incrementais, ou seja, calculamos novamente os atrasos e as receitas, mas com base apenas nos dados filtrados. Os resultados são novamente mostrados na consola e em vez de darmos “overwrite” na coleção do MongoDB, fazemos “append” de forma a juntar aos dados de dias anteriores, sem os subsituir."

$anchor.InsertParagraphAfter()
$anchor.Collapse(0)
$anchor.Move(1, 1) | Out-Null
$anchor.InsertAfter($para1)

$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$anchor.Collapse(0)
$anchor.Move(1, 1) | Out-Null
$anchor.InsertAfter($para2)

$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$anchor.Collapse(0)
$anchor.Move(1, 1) | Out-Null
$anchor.InsertAfter($para3)

$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

Write-Output "done"
